$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(64, 1).Value = 46013
$ws.Cells.Item(64, 2).Value = 137
$ws.Cells.Item(64, 3).Value = 151
$ws.Cells.Item(64, 4).Value = 141

$ws.Cells.Item(64, 1).NumberFormat = $ws.Cells.Item(63, 1).NumberFormat
